# Updated SPI buttons to v1.2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update existing designator lists (new components added to existing groups) ---
$ws.Range("B4").Value = "R1, R2, R3, R4, R5, R6, R7, R8, R9, R10, R11, R12, R13, R14, R15, R16, R17, R18, R19, R20, R21, R22, R23, R24, R25, R26, R27, R28, R29, R30, R31, R32,R36"
$ws.Range("B5").Value = "C1, C2, C3, C4,C8"

# --- Add new BOM rows for v1.2 ---
$ws.Range("A7").Value = "0R"
$ws.Range("B7").Value = "R35,R33"
$ws.Range("C7").Value = "R0603"
$ws.Range("D7").Value = "C21189"

$ws.Range("A8").Value = "SN74LVC2G14DBV"
$ws.Range("B8").Value = "IC6"
$ws.Range("C8").Value = "SOT23-6"
$ws.Range("D8").Value = "C12401"

$ws.Range("A9").Value = "74LVC1G125DBV"
$ws.Range("B9").Value = "IC5"
$ws.Range("C9").Value = "SOT23-5"
$ws.Range("D9").Value = "C23654"

# New data rows (A:C) use the plain body font/formatting (Arial, vertically centered, no wrap)
$abcRange = $ws.Range("A7:C9")
$abcRange.VerticalAlignment = -4108
$abcRange.WrapText = $false
$abcRange.Font.Name = "Arial"

# Column D on the new rows keeps the sheet's default (unformatted) look, like D4:D6 above them
$dRange = $ws.Range("D7:D9")
$dRange.VerticalAlignment = -4108
$dRange.WrapText = $false
$dRange.Font.Name = "宋体"

# --- Normalize A3 / A6 which previously used a separately-defined (but identical) style ---
$ws.Range("A3").VerticalAlignment = -4108
$ws.Range("A3").WrapText = $false
$ws.Range("A3").Font.Name = "Arial"

$ws.Range("A6").VerticalAlignment = -4108
$ws.Range("A6").WrapText = $false
$ws.Range("A6").Font.Name = "Arial"

# --- Normalize D3 so it shares the same wrap style as D2 instead of its own duplicate ---
$ws.Range("D3").VerticalAlignment = -4108
$ws.Range("D3").WrapText = $true
$ws.Range("D3").Font.Name = "Arial"

# --- Move the active selection (mirrors the saved cursor position in the file) ---
$ws.Range("D7").Select()
